$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# NOTE: Price (D) and Volume (E) columns store plain text in the workbook
# (values like "583.47" or "  -3.45%  "). Excel auto-detects plain numeric
# looking strings ("583.54", "1.00", ...) and would silently convert them
# to real numbers (losing formatting / trailing zeros) if assigned as-is.
# A leading apostrophe forces Excel to keep them as text, matching the
# original inline-string cell type, and the apostrophe itself is not
# stored as part of the cell's text.

$ws.Range("D2").Value = '64.223.21'
$ws.Range("E2").Value = '  -2.19%  '

$ws.Range("D3").Value = '3.485.28'
$ws.Range("E3").Value = '  -2.84%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.13%  '

$ws.Range("D5").Value = '''583.54'
$ws.Range("E5").Value = '  -3.25%  '

$ws.Range("D6").Value = '''131.52'
$ws.Range("E6").Value = '  -3.75%  '

$ws.Range("D7").Value = '3.485.88'
$ws.Range("E7").Value = '  -2.85%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("E9").Value = '  -1.68%  '

$ws.Range("E10").Value = '  -0.95%  '

$ws.Range("E11").Value = '  -0.04%  '

$ws.Range("E12").Value = '  -1.48%  '

$ws.Range("D13").Value = '4.069.73'
$ws.Range("E13").Value = '  -3.25%  '

$ws.Range("D14").Value = '''27.76'
$ws.Range("E14").Value = '  -0.80%  '

$ws.Range("E15").Value = '  -4.48%  '

$ws.Range("E16").Value = '  +0.27%  '

$ws.Range("D17").Value = '3.479.31'
$ws.Range("E17").Value = '  -3.20%  '

$ws.Range("D18").Value = '64.241.32'
$ws.Range("E18").Value = '  -2.33%  '

$ws.Range("D19").Value = '''9.94'
$ws.Range("E19").Value = '  -0.89%  '

$ws.Range("D20").Value = '''14.25'
$ws.Range("E20").Value = '  -2.53%  '

$ws.Range("E21").Value = '  -3.98%  '

$ws.Range("D22").Value = '''391.45'
$ws.Range("E22").Value = '  -1.17%  '

$ws.Range("E23").Value = '  -2.36%  '

$ws.Range("D24").Value = '3.624.87'
$ws.Range("E24").Value = '  -2.93%  '

$ws.Range("D25").Value = '''73.10'
$ws.Range("E25").Value = '  -1.64%  '

$ws.Range("D26").Value = '''1.00'
$ws.Range("E26").Value = '  +0.02%  '

$ws.Range("E27").Value = '  -8.14%  '

$ws.Range("E28").Value = '  -6.22%  '

$ws.Range("D29").Value = '''7.48'
$ws.Range("E29").Value = '  -8.58%  '

$ws.Range("E30").Value = '  -0.55%  '

$ws.Range("D31").Value = '''2.24'
$ws.Range("E31").Value = '  -7.07%  '

$ws.Range("E32").Value = '  -4.97%  '

$ws.Range("D33").Value = '3.482.89'
$ws.Range("E33").Value = '  -2.96%  '

$ws.Range("E35").Value = '  -2.60%  '

$ws.Range("E36").Value = '  -2.10%  '

$ws.Range("D37").Value = '''5.27'
$ws.Range("E37").Value = '  -2.23%  '

$ws.Range("B38").Value = 'ImmutableX'
$ws.Range("C38").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D38").Value = '''1.57'
$ws.Range("E38").Value = '  -2.23%  '

$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '''6.97'
$ws.Range("E39").Value = '  -1.55%  '

$ws.Range("D40").Value = '''170.04'
$ws.Range("E40").Value = '  +0.54%  '

$ws.Range("D41").Value = '''0.0805'
$ws.Range("E41").Value = '  -3.71%  '

$ws.Range("E42").Value = '  -3.55%  '

$ws.Range("B43").Value = 'EnergySwap'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D43").Value = '''25.64'
$ws.Range("E43").Value = '  -3.45%  '

$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").Value = '''0.998'
$ws.Range("E44").Value = '  -0.25%  '

$ws.Range("D45").Value = '''41.68'
$ws.Range("E45").Value = '  -3.71%  '

$ws.Range("E46").Value = '  -5.78%  '

$ws.Range("E47").Value = '  -4.13%  '

$ws.Range("E48").Value = '  -3.63%  '

$ws.Range("D49").Value = '''6.86'
$ws.Range("E49").Value = '  -2.57%  '

$ws.Range("D50").Value = '2.422.83'

$ws.Range("D51").Value = '''0.0266'
$ws.Range("E51").Value = '  -1.88%  '
